$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 10429565
$ws.Range("I62").Value = 11124203
$ws.Range("K62").Value = 11124203
$ws.Range("M62").Value = -11123579
$ws.Range("H65").Value = 10429565
$ws.Range("I65").Value = 11124203
$ws.Range("K65").Value = 55621015
$ws.Range("M65").Value = -55617895
$ws.Range("H69").Value = 13999.2
$ws.Range("I69").Value = 12499.25
$ws.Range("K69").Value = 37497.75
$ws.Range("M69").Value = -36623.75
$ws.Range("H72").Value = 13999.2
$ws.Range("I72").Value = 12499.25
$ws.Range("K72").Value = 112493.25
$ws.Range("M72").Value = -108125.25
$ws.Range("H82").Value = 17999.5
$ws.Range("I82").Value = 17999.5
$ws.Range("K82").Value = 53998.5
$ws.Range("M82").Value = -53592.5
$ws.Range("H85").Value = 17999.5
$ws.Range("I85").Value = 17999.5
$ws.Range("K85").Value = 53998.5
$ws.Range("M85").Value = -52594.5
$ws.Range("H88").Value = 2914.9443
$ws.Range("J88").Value = 3139
$ws.Range("L88").Value = 3139
$ws.Range("N88").Value = -3951
$ws.Range("H91").Value = 2914.9443
$ws.Range("J91").Value = 3139
$ws.Range("L91").Value = 3139
$ws.Range("N91").Value = -5947
$ws.Range("H98").Value = 4655.4165
$ws.Range("I98").Value = 4655.4165
$ws.Range("K98").Value = 4655.4165
$ws.Range("M98").Value = -3157.4165
$ws.Range("H122").Value = 4655.4165
$ws.Range("I122").Value = 4655.4165
$ws.Range("K122").Value = 13966.2495
$ws.Range("M122").Value = -11516.2495
$ws.Range("H137").Value = 5723.2856
$ws.Range("I137").Value = 7099.778
$ws.Range("K137").Value = 21299.334
$ws.Range("M137").Value = -18749.334
$ws.Range("H138").Value = 7065.171
$ws.Range("I138").Value = 4723.125
$ws.Range("J138").Value = 7632.9395
$ws.Range("K138").Value = 14169.375
$ws.Range("L138").Value = 22898.8185
$ws.Range("M138").Value = -9029.375
$ws.Range("N138").Value = -33178.8185

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 18560180
$ws.Range("I32").Value = 0
$ws.Range("J32").Value = 18560180
$ws.Range("K32").Value = 0
$ws.Range("L32").Value = 18560180
$ws.Range("M32").ClearContents()
$ws.Range("N32").Value = -18560754
$ws.Range("H61").Value = 10633.643
$ws.Range("I61").Value = 11580.272
$ws.Range("J61").Value = 7162.6665
$ws.Range("K61").Value = 11580.272
$ws.Range("L61").Value = 7162.6665
$ws.Range("M61").Value = -11368.272
$ws.Range("N61").Value = -7586.6665
$ws.Range("H74").Value = 8028.375
$ws.Range("I74").Value = 10033.333
$ws.Range("J74").Value = 2013.5
$ws.Range("K74").Value = 10033.333
$ws.Range("L74").Value = 2013.5
$ws.Range("M74").Value = -9159.333000000001
$ws.Range("N74").Value = -3761.5
$ws.Range("H77").Value = 8028.375
$ws.Range("I77").Value = 10033.333
$ws.Range("J77").Value = 2013.5
$ws.Range("K77").Value = 50166.665
$ws.Range("L77").Value = 10067.5
$ws.Range("M77").Value = -45798.665
$ws.Range("N77").Value = -18803.5
$ws.Range("H80").Value = 26050
$ws.Range("I80").Value = 26050
$ws.Range("K80").Value = 26050
$ws.Range("M80").Value = -25052
$ws.Range("H83").Value = 26050
$ws.Range("I83").Value = 26050
$ws.Range("K83").Value = 78150
$ws.Range("M83").Value = -73158
$ws.Range("H110").Value = 4319.8237
$ws.Range("I110").Value = 1857.0714
$ws.Range("K110").Value = 1857.0714
$ws.Range("M110").Value = 187.9286
$ws.Range("H132").Value = 753119.8
$ws.Range("I132").Value = 971398.2
$ws.Range("J132").Value = 69180.92999999999
$ws.Range("K132").Value = 2914194.6
$ws.Range("L132").Value = 207542.79
$ws.Range("M132").Value = -2911664.6
$ws.Range("N132").Value = -212602.79
$ws.Range("H133").Value = 126298.336
$ws.Range("J133").Value = 126298.336
$ws.Range("L133").Value = 126298.336
$ws.Range("N133").Value = -131358.336
$ws.Range("H135").Value = 100000
$ws.Range("J135").Value = 100000
$ws.Range("L135").Value = 100000
$ws.Range("M135").Value = -110140
$ws.Range("H136").Value = 10633.643
$ws.Range("I136").Value = 11580.272
$ws.Range("J136").Value = 7162.6665
$ws.Range("K136").Value = 34740.81600000001
$ws.Range("L136").Value = 21487.9995
$ws.Range("M136").Value = -32190.81600000001
$ws.Range("N136").Value = -26587.9995
$ws.Range("H139").Value = 76340.25
$ws.Range("J139").Value = 76340.25
$ws.Range("L139").Value = 76340.25
$ws.Range("N139").Value = -86620.25

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H103").Value = 16999
$ws.Range("J103").Value = 16999
$ws.Range("L103").Value = 16999
$ws.Range("N103").Value = -19343

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 142.4375
$ws.Range("I7").Value = 98.166664
$ws.Range("J7").Value = 275.25
$ws.Range("K7").Value = 98.166664
$ws.Range("L7").Value = 275.25
$ws.Range("M7").Value = 14.833336
$ws.Range("N7").Value = -501.25
$ws.Range("H22").Value = 1051148.8
$ws.Range("I22").Value = 1623849.1
$ws.Range("J22").Value = 1198
$ws.Range("K22").Value = 1623849.1
$ws.Range("L22").Value = 1198
$ws.Range("M22").Value = -1623499.1
$ws.Range("N22").Value = -1898
$ws.Range("H31").Value = 26478.834
$ws.Range("I31").Value = 33291
$ws.Range("J31").Value = 19666.666
$ws.Range("K31").Value = 33291
$ws.Range("L31").Value = 19666.666
$ws.Range("M31").Value = -32996
$ws.Range("N31").Value = -20256.666
$ws.Range("H34").Value = 26478.834
$ws.Range("I34").Value = 33291
$ws.Range("J34").Value = 19666.666
$ws.Range("K34").Value = 33291
$ws.Range("L34").Value = 19666.666
$ws.Range("M34").Value = -33089
$ws.Range("N34").Value = -20070.666
$ws.Range("H58").Value = 11620.375
$ws.Range("I58").Value = 8998.299999999999
$ws.Range("J58").Value = 15990.5
$ws.Range("K58").Value = 8998.299999999999
$ws.Range("L58").Value = 15990.5
$ws.Range("M58").Value = -8795.299999999999
$ws.Range("N58").Value = -16396.5
$ws.Range("H60").Value = 37868.332
$ws.Range("J60").Value = 55552.5
$ws.Range("L60").Value = 55552.5
$ws.Range("N60").Value = -56574.5
$ws.Range("H133").Value = 94994.5
$ws.Range("J133").Value = 94994.5
$ws.Range("L133").Value = 94994.5
$ws.Range("N133").Value = -100054.5
$ws.Range("H136").Value = 11620.375
$ws.Range("I136").Value = 8998.299999999999
$ws.Range("J136").Value = 15990.5
$ws.Range("K136").Value = 26994.9
$ws.Range("L136").Value = 47971.5
$ws.Range("M136").Value = -24444.9
$ws.Range("N136").Value = -53071.5

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H106").Value = 29999.834
$ws.Range("J106").Value = 29999.834
$ws.Range("L106").Value = 89999.50199999999
$ws.Range("N106").Value = -91891.50199999999
$ws.Range("H117").Value = 3076.0667
$ws.Range("I117").Value = 302.14285
$ws.Range("K117").Value = 906.4285500000001
$ws.Range("M117").Value = 2535.57145
$ws.Range("H128").Value = 30000
$ws.Range("I128").Value = 30000
$ws.Range("K128").Value = 90000
$ws.Range("M128").Value = -85020
$ws.Range("H131").Value = 6243
$ws.Range("J131").Value = 3257.75
$ws.Range("L131").Value = 9773.25
$ws.Range("N131").Value = -19853.25

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 4650.8184
$ws.Range("I102").Value = 3437
$ws.Range("J102").Value = 7442.6
$ws.Range("K102").Value = 3437
$ws.Range("L102").Value = 7442.6
$ws.Range("M102").Value = -1815
$ws.Range("N102").Value = -10686.6
$ws.Range("H132").Value = 7128.3
$ws.Range("I132").Value = 6756.0303
$ws.Range("J132").Value = 7850.9414
$ws.Range("K132").Value = 20268.0909
$ws.Range("L132").Value = 23552.8242
$ws.Range("M132").Value = -17738.0909
$ws.Range("N132").Value = -28612.8242

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 41668016
$ws.Range("J46").Value = 83335000
$ws.Range("L46").Value = 83335000
$ws.Range("N46").Value = -83335376
$ws.Range("H55").Value = 1419.0834
$ws.Range("I55").Value = 404.83334
$ws.Range("J55").Value = 1757.1666
$ws.Range("K55").Value = 404.83334
$ws.Range("L55").Value = 1757.1666
$ws.Range("M55").Value = -231.83334
$ws.Range("N55").Value = -2103.1666
$ws.Range("H68").Value = 2141.25
$ws.Range("I68").Value = 2199.5454
$ws.Range("K68").Value = 2199.5454
$ws.Range("M68").Value = -1450.5454
$ws.Range("H71").Value = 2141.25
$ws.Range("I71").Value = 2199.5454
$ws.Range("K71").Value = 10997.727
$ws.Range("M71").Value = -7253.726999999999
$ws.Range("H132").Value = 5398.0435
$ws.Range("I132").Value = 3848.2917
$ws.Range("K132").Value = 11544.8751
$ws.Range("M132").Value = -9014.875100000001
$ws.Range("H136").Value = 9811.959999999999
$ws.Range("J136").Value = 8200.200000000001
$ws.Range("L136").Value = 24600.6
$ws.Range("N136").Value = -29700.6

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 8262.161
$ws.Range("I132").Value = 7066.1177
$ws.Range("J132").Value = 20461.8
$ws.Range("K132").Value = 21198.3531
$ws.Range("L132").Value = 61385.39999999999
$ws.Range("M132").Value = -18668.3531
$ws.Range("N132").Value = -66445.39999999999
$ws.Range("H136").Value = 9101199
$ws.Range("I136").Value = 11913531
$ws.Range("K136").Value = 35740593
$ws.Range("M136").Value = -35740593
